# Reorganize feature files and step definitions with improved naming and test data structure
# Append new test result rows (534-542) to the existing ScenarioResults sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Country with Parameter", "PASSED", "chrome"),
    @("Country with Parameter", "PASSED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "FAILED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "FAILED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "FAILED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "FAILED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "FAILED", "chrome"),
    @("Create and Delete CitizenShip From Excel", "PASSED", "chrome")
)

$startRow = 534
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
